# Update the EC (Estado de Cuenta) data table on Hoja1.
# Previous worker/period rows are replaced with a new set of rows:
# JORGE DAVID PETRO MONTES (doc 1133749254) for periods 1911..1903 (rows 17-25)
# JORGE LUIS CARBONELL FLOREZ (doc 1046404745) for periods 1911..1905 (rows 26-32)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Doc type (column B) stays "CC" for every row already - not touched.

# --- JORGE DAVID PETRO MONTES block: rows 17-25 ---
$petroRows = 17..25
$petroPeriods = @("1911","1910","1909","1908","1907","1906","1905","1904","1903")
$petroMora    = @(26500,33125,33125,33125,33125,33125,33125,33125,27604)

for ($i = 0; $i -lt $petroRows.Length; $i++) {
    $r = $petroRows[$i]
    $ws.Range("C$r").Value = "1133749254"
    $ws.Range("D$r").Value = "JORGE DAVID PETRO MONTES"
    $ws.Range("E$r").Value = $petroPeriods[$i]
    $ws.Range("F$r").Value = $petroMora[$i]
}

# --- JORGE LUIS CARBONELL FLOREZ block: rows 26-32 ---
$carbonellRows = 26..32
$carbonellPeriods = @("1911","1910","1909","1908","1907","1906","1905")
$carbonellMora    = @(26500,33125,33125,33125,33125,33125,33125)

for ($i = 0; $i -lt $carbonellRows.Length; $i++) {
    $r = $carbonellRows[$i]
    $ws.Range("C$r").Value = "1046404745"
    $ws.Range("D$r").Value = "JORGE LUIS CARBONELL FLOREZ"
    $ws.Range("E$r").Value = $carbonellPeriods[$i]
    $ws.Range("F$r").Value = $carbonellMora[$i]
}
